# Apply repulled data values to column F ("dSF") on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -5
    3  = 0
    4  = -4
    5  = -3
    8  = 0
    10 = -3
    12 = -1
    15 = 1
    17 = -3
    20 = 0
    22 = -1
    27 = 1
    28 = 2
    33 = -1
    35 = -4
    37 = 6
    43 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
